$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 560.93335
$ws.Range("I39").Value = 60.5
$ws.Range("J39").Value = 1561.8
$ws.Range("K39").Value = 181.5
$ws.Range("L39").Value = 4685.4
$ws.Range("M39").Value = 114.5
$ws.Range("N39").Value = -5277.4
$ws.Range("H76").Value = 3150
$ws.Range("I76").Value = 2695
$ws.Range("K76").Value = 2695
$ws.Range("M76").Value = -2380
$ws.Range("H79").Value = 3150
$ws.Range("I79").Value = 2695
$ws.Range("K79").Value = 2695
$ws.Range("M79").Value = -1603
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1990.0526
$ws.Range("I2").Value = 1423.4445
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 1423.4445
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -1310.4445
$ws.Range("N2").Value = -2726
$ws.Range("H86").Value = 40000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 40000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 40000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -42372
$ws.Range("H89").Value = 40000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 40000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 120000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -131856
$ws.Range("H116").Value = 1990.0526
$ws.Range("I116").Value = 1423.4445
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 1423.4445
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 870.5554999999999
$ws.Range("N116").Value = -7088
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1990.0526
$ws.Range("I3").Value = 1423.4445
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 1423.4445
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -1309.4445
$ws.Range("N3").Value = -2728
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 24683.715
$ws.Range("J54").Value = 24683.715
$ws.Range("L54").Value = 24683.715
$ws.Range("N54").Value = -25999.715
$ws.Range("H99").Value = 3375
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 5750
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 5750
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -8746
$ws.Range("H103").Value = 35000
$ws.Range("I103").Value = 25000
$ws.Range("J103").Value = 40000
$ws.Range("K103").Value = 25000
$ws.Range("L103").Value = 40000
$ws.Range("M103").Value = -23828
$ws.Range("N103").Value = -42344
$ws.Range("H105").Value = 4144.2856
$ws.Range("I105").Value = 3502.5
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3502.5
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1755.5
$ws.Range("N105").Value = -8494
$ws.Range("H107").Value = 1139.6
$ws.Range("I107").Value = 481.75
$ws.Range("K107").Value = 481.75
$ws.Range("M107").Value = 1438.25
$ws.Range("H126").Value = 3375
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -22190
$ws.Range("H134").Value = 4525.3335
$ws.Range("I134").Value = 900
$ws.Range("J134").Value = 4978.5
$ws.Range("K134").Value = 2700
$ws.Range("L134").Value = 14935.5
$ws.Range("M134").Value = -165
$ws.Range("N134").Value = -20005.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5824.3125
$ws.Range("I3").Value = 2246.5557
$ws.Range("J3").Value = 10424.286
$ws.Range("K3").Value = 6739.6671
$ws.Range("L3").Value = 31272.858
$ws.Range("M3").Value = -6627.6671
$ws.Range("N3").Value = -31496.858
$ws.Range("H55").Value = 1612.5
$ws.Range("I55").Value = 580
$ws.Range("J55").Value = 3333.3333
$ws.Range("K55").Value = 1740
$ws.Range("L55").Value = 9999.999899999999
$ws.Range("M55").Value = -1563
$ws.Range("N55").Value = -10353.9999
$ws.Range("H131").Value = 1186.3821
$ws.Range("J131").Value = 1013.1
$ws.Range("L131").Value = 3039.3
$ws.Range("N131").Value = -13119.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 7000
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -6730
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 7000
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -6064
$ws.Range("N73").Value = -8872
$ws.Range("H132").Value = 4052.7666
$ws.Range("I132").Value = 4363.2856
$ws.Range("K132").Value = 13089.8568
$ws.Range("M132").Value = -10559.8568
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1411.5385
$ws.Range("I46").Value = 979.5454999999999
$ws.Range("K46").Value = 979.5454999999999
$ws.Range("M46").Value = -791.5454999999999
$ws.Range("H55").Value = 1139.3077
$ws.Range("I55").Value = 151.5
$ws.Range("J55").Value = 1986
$ws.Range("K55").Value = 151.5
$ws.Range("L55").Value = 1986
$ws.Range("M55").Value = 21.5
$ws.Range("N55").Value = -2332
$ws.Range("H132").Value = 3554.762
$ws.Range("I132").Value = 2387.5
$ws.Range("K132").Value = 7162.5
$ws.Range("M132").Value = -4632.5
$ws.Range("H134").Value = 23000
$ws.Range("J134").Value = 23000
$ws.Range("L134").Value = 23000
$ws.Range("N134").Value = -33140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1763.7142
$ws.Range("I81").Value = 1586.5
$ws.Range("K81").Value = 3173
$ws.Range("M81").Value = -2112
$ws.Range("H84").Value = 1763.7142
$ws.Range("I84").Value = 1586.5
$ws.Range("K84").Value = 15865
$ws.Range("M84").Value = -10561
$ws.Range("H126").Value = 3450561.5
$ws.Range("I126").Value = 1626.091
$ws.Range("J126").Value = 14290073
$ws.Range("K126").Value = 4878.272999999999
$ws.Range("L126").Value = 42870219
$ws.Range("M126").Value = -2408.272999999999
$ws.Range("N126").Value = -42875159
$ws.Range("H132").Value = 4131.04
$ws.Range("I132").Value = 1679.738
$ws.Range("K132").Value = 5039.214
$ws.Range("M132").Value = -2509.214
